# chore(results): Auto-update draw results on excel 2025-12-11T17:52:25Z
# Appends the latest Pick 4 draw result as a new row at the bottom of the
# results table (row 86), mirroring the existing rows' layout/typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 86

# Column A: draw date. Written with a leading apostrophe so Excel keeps it
# as literal text (otherwise "2025-12-11" would be auto-parsed as a date
# serial), then the quote-prefix style is cleared back to Normal so no new
# cell style gets introduced, matching the rest of the column.
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.Value = "'2025-12-11"
$cellA.Style = "Normal"

# Column B: game name - plain text, no ambiguity.
$ws.Cells.Item($newRow, 2).Value = "Pick 4"

# Column C: phase code. Also quote-prefixed so the numeric-looking string
# "251211" is kept as text rather than becoming the number 251211.
$cellC = $ws.Cells.Item($newRow, 3)
$cellC.Value = "'251211"
$cellC.Style = "Normal"

# Column D: draw result digits - plain text, no ambiguity.
$ws.Cells.Item($newRow, 4).Value = "4-1-6-4"

# Column E: ISO-8601 insertion timestamp - plain text, no ambiguity.
$ws.Cells.Item($newRow, 5).Value = "2025-12-11T21:52:25.621+04:00"
